# update camas h. SE34
# Append the new weekly rows (epidemiological week 34: 2021-08-22 .. 2021-08-28)
# to the bottom of the "disponibilidad camas hospitalarias" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the date-cell format (style) from the last existing row down
#        onto the new rows before writing values into them -----------------
$ws.Range("A418").Copy($ws.Range("A419:A425"))

# --- 2. Write the new data -------------------------------------------------
$dates = @(44430, 44431, 44432, 44433, 44434, 44435, 44436)
$bvals = @(1, 1, 0.94871794871794868, 0.94871794871794868, 0.94871794871794868, 0.97435897435897434, 0.97435897435897434)
$cvals = @(0.16483516483516483, 0.1650943396226415, 0.16981132075471697, 0.15748031496062992, 0.17637795275590551, 0.17535545023696683, 0.17535545023696683)
$dvals = @(0.088, 0.10843373493975904, 0.11244979919678715, 0.11693548387096774, 0.10887096774193548, 0.10975609756097561, 0.10975609756097561)
$evals = @(0.64, 0.6, 0.52, 0.4, 0.44, 0.36, 0.36)

for ($i = 0; $i -lt 7; $i++) {
    $r = 419 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
    $ws.Cells.Item($r, 5).Value = $evals[$i]
}

# --- 3. Scroll / select like the author left the sheet after editing ------
$win = $excel.ActiveWindow
$win.ScrollRow = 398
$win.ScrollColumn = 1
$null = $ws.Range("A418:A425").Select()

# --- 4. Reproduce the (orphaned) conditional-format styles ----------------
#        The author applied & removed a "Duplicate Values" highlight twice,
#        which leaves two unused dxf entries behind in styles.xml without
#        any conditionalFormatting rule remaining on the sheet.
for ($n = 0; $n -lt 2; $n++) {
    $fc = $ws.Range("A2:A425").FormatConditions()
    $cond = $fc.AddUniqueValues()
    $cond.Font.Color = 26012
    $cond.Interior.Color = 10284031
    $cond.Delete()
}

Write-Output "done"
